$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..3) {
    $ws.Range("D$row").Value = -0.061

    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = -2.3
    $ws.Range("L$row").Value = -0.1133004926108374
    $ws.Range("M$row").Value = 7
    $ws.Range("N$row").Value = 0.0957592339261286
    $ws.Range("O$row").Value = -3.043478260869565
    $ws.Range("P$row").Value = 7
    $ws.Range("Q$row").Value = 0.0957592339261286
    $ws.Range("R$row").Value = -3.043478260869565

    $ws.Range("U$row").Value = 55.5
    $ws.Range("V$row").Value = 0.7592339261285911
    $ws.Range("W$row").Value = -0.0224609375
    $ws.Range("X$row").Value = 0.08752619807710092
    $ws.Range("Y$row").Value = -0.1099871355771009
    $ws.Range("Z$row").Value = 0.4126016260162602
    $ws.Range("AA$row").Value = 0
    $ws.Range("AB$row").Value = 0.08750317121269791
    $ws.Range("AC$row").Value = -0.08750317121269791
    $ws.Range("AD$row").Value = 0.3
    $ws.Range("AE$row").Value = 0
    $ws.Range("AF$row").Value = 0.3
    $ws.Range("AG$row").Value = -55.2
    $ws.Range("AH$row").Value = 0.004087193460490464
    $ws.Range("AI$row").Value = 0.0027124773960217
    $ws.Range("AJ$row").Value = -3.083798882681566
    $ws.Range("AK$row").Value = -1.001814882032668

    $ws.Range("AN$row").ClearContents()
    $ws.Range("AP$row").ClearContents()
}
